# Apply the diff: insert two new price-report rows (weekly update) at rows 33-34,
# shifting the existing data rows 33-112 down to 35-114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 33 - existing rows 33-112 shift down to 35-114.
$ws.Range("A33:A34").EntireRow.Insert()

# New row 33: Ají / Inferno / Primera
$ws.Range("A33").Value = 8
$ws.Range("B33").Value = "Terminal La Palmera de La Serena"
$ws.Range("C33").Value = "Coquimbo"
$ws.Range("D33").Value = "2021-10-29"
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 100112021
$ws.Range("G33").Value = "Ají"
$ws.Range("H33").Value = "Inferno"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 560
$ws.Range("K33").Value = 33000
$ws.Range("L33").Value = 34000
$ws.Range("M33").Value = 33500
$ws.Range("N33").Value = "$/caja 12 kilos"
$ws.Range("O33").Value = "Región de Arica y Parinacota"
$ws.Range("P33").Value = 2792
$ws.Range("Q33").Value = 12
$ws.Range("R33").Value = "Hortaliza"

# New row 34: Ají / Inferno / Segunda
$ws.Range("A34").Value = 8
$ws.Range("B34").Value = "Terminal La Palmera de La Serena"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = "2021-10-29"
$ws.Range("E34").Value = 4
$ws.Range("F34").Value = 100112021
$ws.Range("G34").Value = "Ají"
$ws.Range("H34").Value = "Inferno"
$ws.Range("I34").Value = "Segunda"
$ws.Range("J34").Value = 360
$ws.Range("K34").Value = 24000
$ws.Range("L34").Value = 25000
$ws.Range("M34").Value = 24500
$ws.Range("N34").Value = "$/caja 12 kilos"
$ws.Range("O34").Value = "Región de Arica y Parinacota"
$ws.Range("P34").Value = 2042
$ws.Range("Q34").Value = 12
$ws.Range("R34").Value = "Hortaliza"
